$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.092.82"
$ws.Range("E2").Value = "  +4.82%  "

$ws.Range("D3").Value = "3.243.61"
$ws.Range("E3").Value = "  +2.15%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "394.85"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.75"
$ws.Range("E6").Value = "  -1.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.587"
$ws.Range("E7").Value = "  +6.79%  "

$ws.Range("D8").Value = "3.237.99"
$ws.Range("E8").Value = "  +2.08%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.623"
$ws.Range("E10").Value = "  +0.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "39.05"
$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0989"
$ws.Range("E12").Value = "  +11.90%  "

$ws.Range("E13").Value = "  +1.66%  "

$ws.Range("D14").Value = "3.763.80"
$ws.Range("E14").Value = "  +2.40%  "

$ws.Range("E15").Value = "  +2.21%  "

$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("D17").Value = "3.232.94"
$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("E18").Value = "  -2.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.76"
$ws.Range("E19").Value = "  +2.46%  "

$ws.Range("D20").Value = "56.962.96"
$ws.Range("E20").Value = "  +4.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.32"
$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000111"
$ws.Range("E22").Value = "  +12.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.90"
$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "295.65"
$ws.Range("E24").Value = "  +8.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.16"
$ws.Range("E25").Value = "  +3.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.16"
$ws.Range("E26").Value = "  -2.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.95"
$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.57"
$ws.Range("E28").Value = "  -5.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.23"
$ws.Range("E29").Value = "  -1.79%  "

$ws.Range("E30").Value = "  -1.33%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.30"
$ws.Range("E32").Value = "  +2.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("E33").Value = "  -3.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "39.14"
$ws.Range("E34").Value = "  +6.17%  "

$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("E36").Value = "  +2.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.46"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("E39").Value = "  -4.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  +2.75%  "

$ws.Range("E41").Value = "  +4.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "134.67"
$ws.Range("E42").Value = "  +3.45%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.88"
$ws.Range("E43").Value = "  -2.07%  "

$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.00"
$ws.Range("E44").Value = "  -1.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.93"
$ws.Range("E45").Value = "  -4.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.279"
$ws.Range("E46").Value = "  -3.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.07"
$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("E48").Value = "  +3.37%  "

$ws.Range("D49").Value = "2.154.08"
$ws.Range("E49").Value = "  +3.02%  "

$ws.Range("E50").Value = "  -5.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.94"
$ws.Range("E51").Value = "  +15.48%  "
